$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Unprotect()

# Update the confidential disclaimer text (cell A16) - date changed from 2021-04-21 to 2021-04-22
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-22 for illustrative purposes only and are subject to change."

# Update the performance figures in columns D (return) and E (excess return) for rows 2-13
$ws.Range("D2").Value = 0.03083680604270389
$ws.Range("E2").Value = -0.004787961696306242

$ws.Range("D3").Value = 0.02420970582770558
$ws.Range("E3").Value = -0.007131354957441993

$ws.Range("D4").Value = 0.05174250065409539
$ws.Range("E4").Value = -0.007241298761971438

$ws.Range("D5").Value = 0.1390994291358736
$ws.Range("E5").Value = -0.002898084044437144

$ws.Range("D6").Value = 0.02851820375924739
$ws.Range("E6").Value = -0.01320901320901324

$ws.Range("D7").Value = 0.1209934287111856
$ws.Range("E7").Value = -0.007711442786069833

$ws.Range("D8").Value = 0.1002578936990706
$ws.Range("E8").Value = -0.003548085901027176

$ws.Range("D9").Value = 0.02798639576373695
$ws.Range("E9").Value = -0.0161952062189592

$ws.Range("D10").Value = 0.1195383034471775
$ws.Range("E10").Value = -0.01050753370340995

$ws.Range("D11").Value = 0.2547926707188982
$ws.Range("E11").Value = -0.01007334099142876

$ws.Range("D12").Value = 0.1020246622403052
$ws.Range("E12").Value = -0.005853658536585371

$ws.Range("D13").Value = 0.9999999999999999
$ws.Range("E13").Value = -0.007636682238560955

$ws.Protect("", $true, $true, $true)

$wb.Save()
